# Update data for Washington, D.C. (p_hp_c_inv) across scenario sheets

$wb = $excel.ActiveWorkbook

# 2025: new base value (plain number, no formula)
$ws2025 = $wb.Worksheets.Item("2025")
$ws2025.Range("A2").Value = 755876

# 2030: formula referencing 2025 base value, with a 20% reduction factor applied
$ws2030 = $wb.Worksheets.Item("2030")
$ws2030.Range("A2").Formula = "='2025'!A2*(1-0.25*0.2)"

# 2035: 40% reduction factor
$ws2035 = $wb.Worksheets.Item("2035")
$ws2035.Range("A2").Formula = "='2025'!A2*(1-0.25*0.4)"

# 2040: 60% reduction factor
$ws2040 = $wb.Worksheets.Item("2040")
$ws2040.Range("A2").Formula = "='2025'!A2*(1-0.25*0.6)"

# 2045: 80% reduction factor
$ws2045 = $wb.Worksheets.Item("2045")
$ws2045.Range("A2").Formula = "='2025'!A2*(1-0.25*0.8)"

# 2050: 100% reduction factor
$ws2050 = $wb.Worksheets.Item("2050")
$ws2050.Range("A2").Formula = "='2025'!A2*(1-0.25*1)"
